$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.07287566666666666
$ws.Range("H2").Value = 0.218627
$ws.Range("I2").Value = 0.6403049419813613
$ws.Range("J2").Value = 0.6403049419813615
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8232033333333333
$ws.Range("N2").Value = 2.46961
$ws.Range("O2").Value = 0.03174646187333627
$ws.Range("P2").Value = 0.03174646187333627
$ws.Range("Q2").Value = 0.05999149171888888
$ws.Range("R2").Value = 0.5399234254699999
$ws.Range("S2").Value = 0.02032741642792008
$ws.Range("T2").Value = 0.02032741642792008

$ws.Range("G3").Value = 0.07287566666666666
$ws.Range("H3").Value = 0.218627
$ws.Range("I3").Value = 0.6403049419813613
$ws.Range("J3").Value = 0.6403049419813615
$ws.Range("O3").Value = 0.4477531724126574
$ws.Range("P3").Value = 0.4477531724126574
$ws.Range("Q3").Value = 0.8461220290334442
$ws.Range("R3").Value = 7.615098261300999
$ws.Range("S3").Value = 0.2866985690836571
$ws.Range("T3").Value = 0.2866985690836571

$ws.Range("G4").Value = 0.07287566666666666
$ws.Range("H4").Value = 0.218627
$ws.Range("I4").Value = 0.6403049419813613
$ws.Range("J4").Value = 0.6403049419813615
$ws.Range("O4").Value = 0.5205003657140063
$ws.Range("P4").Value = 0.5205003657140063
$ws.Range("Q4").Value = 0.9835928647417777
$ws.Range("R4").Value = 8.852335782676001
$ws.Range("S4").Value = 0.3332789564697841
$ws.Range("T4").Value = 0.3332789564697842

$ws.Range("I5").Value = 0.3596950580186386
$ws.Range("J5").Value = 0.3596950580186386
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8232033333333333
$ws.Range("N5").Value = 2.46961
$ws.Range("O5").Value = 0.03174646187333627
$ws.Range("P5").Value = 0.03174646187333627
$ws.Range("Q5").Value = 0.03370057246111111
$ws.Range("R5").Value = 0.30330515215
$ws.Range("S5").Value = 0.01141904544541619
$ws.Range("T5").Value = 0.01141904544541619

$ws.Range("I6").Value = 0.3596950580186386
$ws.Range("J6").Value = 0.3596950580186386
$ws.Range("O6").Value = 0.4477531724126574
$ws.Range("P6").Value = 0.4477531724126574
$ws.Range("R6").Value = 4.277826128345
$ws.Range("S6").Value = 0.1610546033290003
$ws.Range("T6").Value = 0.1610546033290003

$ws.Range("I7").Value = 0.3596950580186386
$ws.Range("J7").Value = 0.3596950580186386
$ws.Range("O7").Value = 0.5205003657140063
$ws.Range("P7").Value = 0.5205003657140063
$ws.Range("R7").Value = 4.97285156522
$ws.Range("S7").Value = 0.1872214092442221
$ws.Range("T7").Value = 0.1872214092442221
